$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.999.34"
$ws.Range("E2").Value = "  +8.16%  "

$ws.Range("D3").Value = "3.142.42"
$ws.Range("E3").Value = "  +5.47%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'590.09"
$ws.Range("E5").Value = "  +4.73%  "

$ws.Range("D6").Value = "'147.30"
$ws.Range("E6").Value = "  +7.52%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.135.74"
$ws.Range("E8").Value = "  +5.40%  "

$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  +3.61%  "

$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  +20.38%  "

$ws.Range("D11").Value = "'5.75"
$ws.Range("E11").Value = "  +8.37%  "

$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +5.01%  "

$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  +11.84%  "

$ws.Range("D14").Value = "'35.99"
$ws.Range("E14").Value = "  +7.11%  "

$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "3.663.04"
$ws.Range("E16").Value = "  +5.53%  "

$ws.Range("D17").Value = "63.942.94"
$ws.Range("E17").Value = "  +8.04%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.16"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.141.47"
$ws.Range("E19").Value = "  +5.51%  "

$ws.Range("D20").Value = "'473.61"
$ws.Range("E20").Value = "  +8.64%  "

$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  +4.15%  "

$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("D23").Value = "'7.56"
$ws.Range("E23").Value = "  +7.41%  "

$ws.Range("D24").Value = "'13.35"
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("D25").Value = "'82.48"

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'8.65"
$ws.Range("E27").Value = "  +12.04%  "

$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  +6.13%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  +10.44%  "

$ws.Range("D32").Value = "'27.03"
$ws.Range("E32").Value = "  +5.03%  "

$ws.Range("E33").Value = "  +5.09%  "

$ws.Range("D34").Value = "0.0₃0877"
$ws.Range("E34").Value = "  +12.72%  "

$ws.Range("E35").Value = "  +17.56%  "

$ws.Range("E36").Value = "  +5.87%  "

$ws.Range("D37").Value = "'3.40"
$ws.Range("E37").Value = "  +21.18%  "

$ws.Range("D38").Value = "'6.14"
$ws.Range("E38").Value = "  +4.35%  "

$ws.Range("D39").Value = "'50.88"
$ws.Range("E39").Value = "  +4.90%  "

$ws.Range("D40").Value = "'445.53"
$ws.Range("E40").Value = "  +12.57%  "

$ws.Range("D41").Value = "'8.75"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0372"
$ws.Range("E42").Value = "  +6.18%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.925.30"
$ws.Range("E43").Value = "  +7.80%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.283"
$ws.Range("E44").Value = "  +13.59%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.112"
$ws.Range("E45").Value = "  +7.65%  "

$ws.Range("D46").Value = "'2.20"
$ws.Range("E46").Value = "  +11.46%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'123.66"
$ws.Range("E48").Value = "  +1.29%  "

$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'34.60"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").Value = "'0.112"
$ws.Range("E50").Value = "  +2.33%  "

$ws.Range("D51").Value = "'24.89"
$ws.Range("E51").Value = "  +7.52%  "
